$wb = $excel.ActiveWorkbook

# "Jaana" sheet gets a new time-tracking entry for the Sprint 2 planning
# meeting (23.02.2023): date, 1 hour, and a description of the meeting.
$wsJaana = $wb.Worksheets.Item("Jaana")

# Copy the date-number-format from the row above (A7) so the new date cell
# picks up the same style as the other date cells in the column.
$wsJaana.Range("A7").Copy()
$wsJaana.Range("A11").PasteSpecial(-4122)

$wsJaana.Range("A11").Value = 44980
$wsJaana.Range("B11").Value = 1
$wsJaana.Range("C11").Value = "Scrum-tiimin palavereja: daily, retro, review ja uuden sprintin planning."

# The "Jarno" sheet's cursor returns to its earlier resting spot (A11) ...
$wsJarno = $wb.Worksheets.Item("Jarno")
$wsJarno.Range("A11").Select()

# ... while "Jaana" becomes the active tab, with the selection parked on
# the next empty row beneath the entry just added.
$wsJaana.Activate()
$wsJaana.Range("B12").Select()
